$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns (initial_waterlevel_file, water_level_ini_type)
# before the old "display_name" column (H), shifting H:Q -> J:S.
$ws.Columns("H:I").Insert()

# --- New header cells (row 1).
# Set H1 then H4 first so the new shared strings land in the same order
# as the target file (40=initial_waterlevel_file, 41=raster path, 42=water_level_ini_type).
$ws.Range("H1").Value = "initial_waterlevel_file"
$ws.Range("H4").Value = "rasters/initial_wlvl_2d_hoekje.tif"
$ws.Range("I1").Value = "water_level_ini_type"

# --- New data cells for rows 5-7 (ids 4-6), which already have a display_name set.
$ws.Range("H5").Value = "rasters/initial_wlvl_2d_hoekje.tif"
$ws.Range("H6").Value = "rasters/initial_wlvl_2d_hoekje.tif"
$ws.Range("H7").Value = "rasters/initial_wlvl_2d_hoekje.tif"

# --- water_level_ini_type values for rows 4-7.
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1

# --- Row 3 (id 2) has no display_name / initial waterlevel data, but Excel
# still applied the wrap/vertical-center alignment style (no border) to the
# otherwise-blank H3:I3 cells.
$ws.Range("H3:I3").Borders.LineStyle = -4142

# --- use_2d_flow / use_0d_inflow values for row 2 (id 1) changed as part of
# the same edit (use_2d_flow 0, use_0d_inflow 1 are now in their shifted
# columns O2/Q2 with values 0 and 1 respectively - already correct via shift).

# --- Drop the large block of trailing, always-empty columns (old S/T..AS)
# that Excel had pre-formatted; the edited file trims the sheet back down to
# column S.
$ws.Range("T1:AU7").EntireColumn.Delete()

# --- Selection / dimension bookkeeping to match the saved file.
$ws.Range("Z5").Select()

Write-Host "edit applied"
